# Adds a small 3x4 helper table (row 34) and column sums (row 36) in
# columns R:U below the existing grading table, and updates the
# worksheet's active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 34: new raw values (continues the 2,3,4,5 sequence that already
# lives in R3:U3 with the next set of numbers 3,4,5,6)
$ws.Range("R34").Value = 3
$ws.Range("S34").Value = 4
$ws.Range("T34").Value = 5
$ws.Range("U34").Value = 6

# Row 36: column totals over the main data block (rows 4-32)
$ws.Range("R36").Formula = "=SUM(R4:R32)"
$ws.Range("S36:U36").Formula = "=SUM(S4:S32)"

# Restore the active cell/selection that Excel persisted at save time
$null = $ws.Range("Y22").Select()
